# edit.ps1 -- applies the "add 2022-Q3 data" change to the workbook.
#
# Summary of the change:
#   1. Insert a new worksheet named "2022-Q3" right after "总计" and before
#      "2022-Q2" (pushing all the quarter sheets down by one position).
#   2. Populate the new sheet with the fund holdings table for 2022-Q3
#      (header row + 28 data rows, columns A-H).
#   3. Insert a new row at the top of the data in the "总计" (totals) sheet
#      with the 2022-Q3 summary (count=28, value=5.09), shifting the
#      existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Fill in the "2022-Q3" fund holdings table.
# ---------------------------------------------------------------------
$headers = @('基金代码','基金名称','基金规模','股票总仓位','仓位占比','持有市值(亿元)','仓位排名')
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$data = @(
    @('0','010659','民生加银质量领先混合A','17.99','89.90','6.47','1.1640','1'),
    @('1','009362','招商丰盈积极配置混合A','11.50','84.77','3.96','0.4554','5'),
    @('2','217009','招商核心价值混合','7.66','79.53','5.31','0.4067','2'),
    @('3','217010','招商大盘蓝筹混合','7.41','78.95','5.37','0.3979','2'),
    @('4','000408','民生加银城镇化混合A','5.81','88.42','6.46','0.3753','1'),
    @('5','010795','民生加银价值发现一年持有期混合A','5.33','91.31','6.48','0.3454','1'),
    @('6','011843','民生加银内核驱动混合A','4.27','89.74','6.49','0.2771','1'),
    @('7','012003','招商价值成长混合A','8.47','86.90','3.15','0.2668','7'),
    @('8','690005','民生加银内需增长混合','3.50','89.95','6.42','0.2247','1'),
    @('9','013559','招商均衡回报混合A','3.62','83.32','5.51','0.1995','3'),
    @('10','011690','招商品质发现混合A','5.99','84.71','3.24','0.1941','6'),
    @('11','009363','招商丰盈积极配置混合C','3.59','84.77','3.96','0.1422','5'),
    @('12','009695','招商成长精选一年定期开放混合A','2.13','90.95','6.66','0.1419','3'),
    @('13','002249','招商境远灵活配置混合','1.53','84.27','6.47','0.0990','2'),
    @('14','012004','招商价值成长混合C','3.13','86.90','3.15','0.0986','7'),
    @('15','005521','华安红利精选混合A','2.43','86.40','2.91','0.0707','9'),
    @('16','010660','民生加银质量领先混合C','0.86','89.90','6.47','0.0556','1'),
    @('17','009696','招商成长精选一年定期开放混合C','0.71','90.95','6.66','0.0473','3'),
    @('18','011727','工银瑞信聚瑞混合A','1.23','38.28','2.19','0.0269','5'),
    @('19','930602','国信价值智选混合','0.49','76.70','5.32','0.0261','5'),
    @('20','010796','民生加银价值发现一年持有期混合C','0.37','91.31','6.48','0.0240','1'),
    @('21','009706','民生加银城镇化混合C','0.34','88.42','6.46','0.0220','1'),
    @('22','013560','招商均衡回报混合C','0.20','83.32','5.51','0.0110','3'),
    @('23','011844','民生加银内核驱动混合C','0.14','89.74','6.49','0.0091','1'),
    @('24','012495','民生加银双核动力混合','0.16','53.77','3.96','0.0063','1'),
    @('25','011691','招商品质发现混合C','0.12','84.71','3.24','0.0039','6'),
    @('26','011728','工银瑞信聚瑞混合C','0.00','38.28','2.19','0','5'),
    @('27','014971','华安红利精选混合C','0.00','86.40','2.91','0','9')
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2

    # Column A: numeric row index (0-based) -- bold/bordered/centered,
    # matching the index-column style used throughout the workbook.
    $aCell = $q3.Cells.Item($excelRow, 1)
    $aCell.Value = [int]$row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # Column B: fund code -- force text so leading zeros survive
    $q3.Cells.Item($excelRow, 2).Value = "'" + $row[1]

    # Column C: fund name (text)
    $q3.Cells.Item($excelRow, 3).Value = "'" + $row[2]

    # Column D: fund scale -- stored as text in the source data
    $q3.Cells.Item($excelRow, 4).Value = "'" + $row[3]

    # Column E: total stock position -- stored as text
    $q3.Cells.Item($excelRow, 5).Value = "'" + $row[4]

    # Column F: position ratio -- stored as text
    $q3.Cells.Item($excelRow, 6).Value = "'" + $row[5]

    # Column G: holding market value -- text, except the two zero rows
    #           which are stored as a real number 0.
    if ($row[6] -eq "0") {
        $q3.Cells.Item($excelRow, 7).Value = 0
    } else {
        $q3.Cells.Item($excelRow, 7).Value = "'" + $row[6]
    }

    # Column H: position rank (real number)
    $q3.Cells.Item($excelRow, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------
# 3) Insert the new summary row into "总计" and fill it in.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalA2 = $totalSheet.Cells.Item(2, 1)
$totalA2.Value = 0
$totalA2.Font.Bold = $true
$totalA2.HorizontalAlignment = -4108
$totalA2.VerticalAlignment = -4160
$totalA2.Borders.LineStyle = 1
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 28
$totalSheet.Cells.Item(2, 4).Value = 5.09

# Column A on the rest of the rows is a 0-based position index, so every
# row below the newly-inserted one needs to be bumped up by one.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Output "2022-Q3 sheet inserted and populated; 总计 updated."
